$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-51 down to 43-52
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new data record
$ws.Cells.Item(42, 1).Value = 10
$ws.Cells.Item(42, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value = "La Araucanía"
$ws.Cells.Item(42, 4).Value = 44637
$ws.Cells.Item(42, 5).Value = 9
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100107
$ws.Cells.Item(42, 8).Value = "Otros"
$ws.Cells.Item(42, 9).Value = 100107011
$ws.Cells.Item(42, 10).Value = "Tuna"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 65
$ws.Cells.Item(42, 14).Value = 15000
$ws.Cells.Item(42, 15).Value = 15000
$ws.Cells.Item(42, 16).Value = 15000
$ws.Cells.Item(42, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(42, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(42, 19).Value = 938
$ws.Cells.Item(42, 20).Value = 16

# Apply the same date number format as column D uses elsewhere
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
